# Add a new "Groups" worksheet at the end of the workbook (after "Trees"),
# populate it with an id / group_name table, and leave it as the active sheet
# — mirroring "added groups model with relationship".

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet ("Trees") so it lands
# at the end of the tab strip and becomes the active sheet, same as Excel
# does when you click "Insert Sheet" on the last tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Groups"

# Header
$ws.Range("A1").Value = "id"

# Data rows first, header text last, so the shared-string table is built in
# the same order ("group1", "group2", "group3", then "group_name").
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "group1"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "group2"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "group3"

$ws.Range("B1").Value = "group_name"

# Leave the selection on the header of the new group-name column.
$ws.Range("B1").Select()
